$wb = $excel.ActiveWorkbook

# Add the two new sheets at the end of the workbook, after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSame = $wb.Worksheets.Add($null, $lastSheet)
$wsSame.Name = "same_elements"

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPartly = $wb.Worksheets.Add($null, $lastSheet2)
$wsPartly.Name = "partly_same"

# --- same_elements sheet data ---
$wsSame.Range("B1").Value = 5
$wsSame.Range("C1").Value = 50
$wsSame.Range("D1").Value = 500
$wsSame.Range("E1").Value = 5000
$wsSame.Range("F1").Value = 50000
$wsSame.Range("G1").Value = 500000

$wsSame.Range("A2").Value = "byte"
$wsSame.Range("B2").Value = 0
$wsSame.Range("C2").Value = 0
$wsSame.Range("D2").Value = 0.0030010000000000002
$wsSame.Range("E2").Value = 0.040002999999999997
$wsSame.Range("F2").Value = 0.44103999999999999
$wsSame.Range("G2").Value = 4.8714490000000001

$wsSame.Range("A3").Value = "int"
$wsSame.Range("B3").Value = 0
$wsSame.Range("C3").Value = 0
$wsSame.Range("D3").Value = 0.0029989999999999999
$wsSame.Range("E3").Value = 0.032002999999999997
$wsSame.Range("F3").Value = 0.43403799999999998
$wsSame.Range("G3").Value = 4.6234080000000004

$wsSame.Range("A4").Value = "string"
$wsSame.Range("B4").Value = 0
$wsSame.Range("C4").Value = 0.00099900000000000001
$wsSame.Range("D4").Value = 0.0030010000000000002
$wsSame.Range("E4").Value = 0.036004000000000001
$wsSame.Range("F4").Value = 0.40003499999999997
$wsSame.Range("G4").Value = 4.6253919999999997

$wsSame.Range("A5").Value = "date"
$wsSame.Range("B5").Value = 0
$wsSame.Range("C5").Value = 0
$wsSame.Range("D5").Value = 0.0030010000000000002
$wsSame.Range("E5").Value = 0.034002999999999999
$wsSame.Range("F5").Value = 0.46004099999999998
$wsSame.Range("G5").Value = 4.9014350000000002

# --- partly_same sheet data ---
$wsPartly.Range("B1").Value = 5
$wsPartly.Range("C1").Value = 50
$wsPartly.Range("D1").Value = 500
$wsPartly.Range("E1").Value = 5000
$wsPartly.Range("F1").Value = 50000
$wsPartly.Range("G1").Value = 500000

$wsPartly.Range("A2").Value = "byte"
$wsPartly.Range("B2").Value = 0
$wsPartly.Range("C2").Value = 0
$wsPartly.Range("D2").Value = 0.0030000000000000001
$wsPartly.Range("E2").Value = 0.031001999999999998
$wsPartly.Range("F2").Value = 0.36903200000000003
$wsPartly.Range("G2").Value = 4.5644039999999997

$wsPartly.Range("A3").Value = "int"
$wsPartly.Range("B3").Value = 0
$wsPartly.Range("C3").Value = 0
$wsPartly.Range("D3").Value = 0.0030000000000000001
$wsPartly.Range("E3").Value = 0.031002999999999999
$wsPartly.Range("F3").Value = 0.37303999999999998
$wsPartly.Range("G3").Value = 4.6514480000000002

$wsPartly.Range("A4").Value = "string"
$wsPartly.Range("B4").Value = 0
$wsPartly.Range("C4").Value = 0
$wsPartly.Range("D4").Value = 0.0030010000000000002
$wsPartly.Range("E4").Value = 0.039003000000000003
$wsPartly.Range("F4").Value = 0.38806499999999999
$wsPartly.Range("G4").Value = 4.9564409999999999

$wsPartly.Range("A5").Value = "date"
$wsPartly.Range("B5").Value = 0
$wsPartly.Range("C5").Value = 0
$wsPartly.Range("D5").Value = 0.0029989999999999999
$wsPartly.Range("E5").Value = 0.036003000000000003
$wsPartly.Range("F5").Value = 0.36999599999999999
$wsPartly.Range("G5").Value = 5.6967059999999998

# Mark same_elements as the active/selected sheet (tabSelected) with selection at M16
$wsSame.Select()
$wsSame.Range("M16").Select()
